# Append SE-run rows 88-108 (user_data.xlsx log rows) to Sheet1.
# Columns A,B are numeric (Date_seconds, ParentCorpID); C-H are free-text log
# fields (rbLocations, cbGroupBizType, txtNumEmployees, txtAssetValue,
# cbBusinessFunctions, Date) that must stay Text even when their content is
# numeric-looking (" 250 ", " 999 ", ...), preserving the leading/trailing
# spaces baked into the source log. A leading "'" forces Text entry the way
# Excel's UI would; re-applying the "Normal" style afterwards clears the
# quote-prefix flag the apostrophe leaves behind, without touching the value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 88
$ws.Range("A88").Value = 1519400664
$ws.Range("B88").Value = 1
$ws.Range("C88").Value = "' Boise "
$ws.Range("C88").Style = "Normal"
$ws.Range("D88").Value = "' - "
$ws.Range("D88").Style = "Normal"
$ws.Range("E88").Value = "' 250 "
$ws.Range("E88").Style = "Normal"
$ws.Range("F88").Value = "' 1000 "
$ws.Range("F88").Style = "Normal"
$ws.Range("G88").Value = "' Clean Room Manufacturing, R&D "
$ws.Range("G88").Style = "Normal"
$ws.Range("H88").Value = "' Fri_Feb_23_10:44:24_EST_2018"
$ws.Range("H88").Style = "Normal"

# Row 89
$ws.Range("A89").Value = 1519403232
$ws.Range("B89").Value = 1
$ws.Range("C89").Value = "' Manassas "
$ws.Range("C89").Style = "Normal"
$ws.Range("D89").Value = "' - "
$ws.Range("D89").Style = "Normal"
$ws.Range("E89").Value = "' 250 "
$ws.Range("E89").Style = "Normal"
$ws.Range("F89").Value = "' 999 "
$ws.Range("F89").Style = "Normal"
$ws.Range("G89").Value = "' Clean Room Manufacturing, R&D "
$ws.Range("G89").Style = "Normal"
$ws.Range("H89").Value = "' Fri_Feb_23_11:27:12_EST_2018"
$ws.Range("H89").Style = "Normal"

# Row 90
$ws.Range("A90").Value = 1519404127
$ws.Range("B90").Value = 4
$ws.Range("C90").Value = "' Indooroopilly_Shopping_Centre "
$ws.Range("C90").Style = "Normal"
$ws.Range("D90").Value = "' - "
$ws.Range("D90").Style = "Normal"
$ws.Range("E90").Value = "' 999 "
$ws.Range("E90").Style = "Normal"
$ws.Range("F90").Value = "' 997 "
$ws.Range("F90").Style = "Normal"
$ws.Range("G90").Value = "' Clean Room Manufacturing, R&D "
$ws.Range("G90").Style = "Normal"
$ws.Range("H90").Value = "' Fri_Feb_23_11:42:07_EST_2018"
$ws.Range("H90").Style = "Normal"

# Row 91
$ws.Range("A91").Value = 1519404349
$ws.Range("B91").Value = 4
$ws.Range("C91").Value = "' Indooroopilly_Shopping_Centre "
$ws.Range("C91").Style = "Normal"
$ws.Range("D91").Value = "' - "
$ws.Range("D91").Style = "Normal"
$ws.Range("E91").Value = "' 999 "
$ws.Range("E91").Style = "Normal"
$ws.Range("F91").Value = "' 9999 "
$ws.Range("F91").Style = "Normal"
$ws.Range("G91").Value = "' Clean Room Manufacturing, R&D "
$ws.Range("G91").Style = "Normal"
$ws.Range("H91").Value = "' Fri_Feb_23_11:45:49_EST_2018"
$ws.Range("H91").Style = "Normal"

# Row 92
$ws.Range("A92").Value = 1519411962
$ws.Range("B92").Value = 4
$ws.Range("C92").Value = "' Coronation_Drive_Office_Park "
$ws.Range("C92").Style = "Normal"
$ws.Range("D92").Value = "' - "
$ws.Range("D92").Style = "Normal"
$ws.Range("E92").Value = "' 250 "
$ws.Range("E92").Style = "Normal"
$ws.Range("F92").Value = "' 999 "
$ws.Range("F92").Style = "Normal"
$ws.Range("G92").Value = "' Clean Room Manufacturing, R&D "
$ws.Range("G92").Style = "Normal"
$ws.Range("H92").Value = "' Fri_Feb_23_13:52:42_EST_2018"
$ws.Range("H92").Style = "Normal"

# Row 93
$ws.Range("A93").Value = 1519411978
$ws.Range("B93").Value = 4
$ws.Range("C93").Value = "' Coronation_Drive_Office_Park "
$ws.Range("C93").Style = "Normal"
$ws.Range("D93").Value = "' - "
$ws.Range("D93").Style = "Normal"
$ws.Range("E93").Value = "' 999 "
$ws.Range("E93").Style = "Normal"
$ws.Range("F93").Value = "' 999 "
$ws.Range("F93").Style = "Normal"
$ws.Range("G93").Value = "' Clean Room Manufacturing, R&D "
$ws.Range("G93").Style = "Normal"
$ws.Range("H93").Value = "' Fri_Feb_23_13:52:58_EST_2018"
$ws.Range("H93").Style = "Normal"

# Row 94
$ws.Range("A94").Value = 1519412355
$ws.Range("B94").Value = 4
$ws.Range("C94").Value = "' Macquarie_Centre "
$ws.Range("C94").Style = "Normal"
$ws.Range("D94").Value = "' - "
$ws.Range("D94").Style = "Normal"
$ws.Range("E94").Value = "' 250 "
$ws.Range("E94").Style = "Normal"
$ws.Range("F94").Value = "' 999 "
$ws.Range("F94").Style = "Normal"
$ws.Range("G94").Value = "' Clean Room Manufacturing, R&D "
$ws.Range("G94").Style = "Normal"
$ws.Range("H94").Value = "' Fri_Feb_23_13:59:15_EST_2018"
$ws.Range("H94").Style = "Normal"

# Row 95
$ws.Range("A95").Value = 1519412485
$ws.Range("B95").Value = 4
$ws.Range("C95").Value = "' Macquarie_Centre "
$ws.Range("C95").Style = "Normal"
$ws.Range("D95").Value = "' - "
$ws.Range("D95").Style = "Normal"
$ws.Range("E95").Value = "' 250 "
$ws.Range("E95").Style = "Normal"
$ws.Range("F95").Value = "' 100 "
$ws.Range("F95").Style = "Normal"
$ws.Range("G95").Value = "' Clean Room Manufacturing, R&D "
$ws.Range("G95").Style = "Normal"
$ws.Range("H95").Value = "' Fri_Feb_23_14:01:25_EST_2018"
$ws.Range("H95").Style = "Normal"

# Row 96
$ws.Range("A96").Value = 1519412536
$ws.Range("B96").Value = 4
$ws.Range("C96").Value = "' Macquarie_Centre "
$ws.Range("C96").Style = "Normal"
$ws.Range("D96").Value = "' - "
$ws.Range("D96").Style = "Normal"
$ws.Range("E96").Value = "' 250 "
$ws.Range("E96").Style = "Normal"
$ws.Range("F96").Value = "' 999 "
$ws.Range("F96").Style = "Normal"
$ws.Range("G96").Value = "' Clean Room Manufacturing, R&D "
$ws.Range("G96").Style = "Normal"
$ws.Range("H96").Value = "' Fri_Feb_23_14:02:16_EST_2018"
$ws.Range("H96").Style = "Normal"

# Row 97
$ws.Range("A97").Value = 1519412681
$ws.Range("B97").Value = 4
$ws.Range("C97").Value = "' Macquarie_Centre "
$ws.Range("C97").Style = "Normal"
$ws.Range("D97").Value = "' - "
$ws.Range("D97").Style = "Normal"
$ws.Range("E97").Value = "' 250 "
$ws.Range("E97").Style = "Normal"
$ws.Range("F97").Value = "' 998 "
$ws.Range("F97").Style = "Normal"
$ws.Range("G97").Value = "' Clean Room Manufacturing, R&D "
$ws.Range("G97").Style = "Normal"
$ws.Range("H97").Value = "' Fri_Feb_23_14:04:41_EST_2018"
$ws.Range("H97").Style = "Normal"

# Row 98
$ws.Range("A98").Value = 1519412850
$ws.Range("B98").Value = 4
$ws.Range("C98").Value = "' Macquarie_Centre "
$ws.Range("C98").Style = "Normal"
$ws.Range("D98").Value = "' - "
$ws.Range("D98").Style = "Normal"
$ws.Range("E98").Value = "' 999 "
$ws.Range("E98").Style = "Normal"
$ws.Range("F98").Value = "' 999 "
$ws.Range("F98").Style = "Normal"
$ws.Range("G98").Value = "' Clean Room Manufacturing, R&D "
$ws.Range("G98").Style = "Normal"
$ws.Range("H98").Value = "' Fri_Feb_23_14:07:30_EST_2018"
$ws.Range("H98").Style = "Normal"

# Row 99
$ws.Range("A99").Value = 1519412851
$ws.Range("B99").Value = 4
$ws.Range("C99").Value = "' Macquarie_Centre "
$ws.Range("C99").Style = "Normal"
$ws.Range("D99").Value = "' - "
$ws.Range("D99").Style = "Normal"
$ws.Range("E99").Value = "' 999 "
$ws.Range("E99").Style = "Normal"
$ws.Range("F99").Value = "' 999 "
$ws.Range("F99").Style = "Normal"
$ws.Range("G99").Value = "' Clean Room Manufacturing, R&D "
$ws.Range("G99").Style = "Normal"
$ws.Range("H99").Value = "' Fri_Feb_23_14:07:31_EST_2018"
$ws.Range("H99").Style = "Normal"

# Row 100
$ws.Range("A100").Value = 1519412868
$ws.Range("B100").Value = 4
$ws.Range("C100").Value = "' AMP_Building "
$ws.Range("C100").Style = "Normal"
$ws.Range("D100").Value = "' - "
$ws.Range("D100").Style = "Normal"
$ws.Range("E100").Value = "' 250 "
$ws.Range("E100").Style = "Normal"
$ws.Range("F100").Value = "' 100 "
$ws.Range("F100").Style = "Normal"
$ws.Range("G100").Value = "' Clean Room Manufacturing, R&D "
$ws.Range("G100").Style = "Normal"
$ws.Range("H100").Value = "' Fri_Feb_23_14:07:48_EST_2018"
$ws.Range("H100").Style = "Normal"

# Row 101
$ws.Range("A101").Value = 1519577983
$ws.Range("B101").Value = 4
$ws.Range("C101").Value = "' AMP_Building "
$ws.Range("C101").Style = "Normal"
$ws.Range("D101").Value = "' - "
$ws.Range("D101").Style = "Normal"
$ws.Range("E101").Value = "' 250 "
$ws.Range("E101").Style = "Normal"
$ws.Range("F101").Value = "' 100 "
$ws.Range("F101").Style = "Normal"
$ws.Range("G101").Value = "' Clean Room Manufacturing, Shipping, R&D "
$ws.Range("G101").Style = "Normal"
$ws.Range("H101").Value = "' Sun_Feb_25_11:59:43_EST_2018"
$ws.Range("H101").Style = "Normal"

# Row 102
$ws.Range("A102").Value = 1519585936
$ws.Range("B102").Value = 4
$ws.Range("C102").Value = "' Indooroopilly_Shopping_Centre "
$ws.Range("C102").Style = "Normal"
$ws.Range("D102").Value = "' - "
$ws.Range("D102").Style = "Normal"
$ws.Range("E102").Value = "' 999 "
$ws.Range("E102").Style = "Normal"
$ws.Range("F102").Value = "' 9999 "
$ws.Range("F102").Style = "Normal"
$ws.Range("G102").Value = "' Clean Room Manufacturing, Inventory Management, R&D "
$ws.Range("G102").Style = "Normal"
$ws.Range("H102").Value = "' Sun_Feb_25_14:12:16_EST_2018"
$ws.Range("H102").Style = "Normal"

# Row 103
$ws.Range("A103").Value = 1519585988
$ws.Range("B103").Value = 4
$ws.Range("C103").Value = "' Indooroopilly_Shopping_Centre "
$ws.Range("C103").Style = "Normal"
$ws.Range("D103").Value = "' - "
$ws.Range("D103").Style = "Normal"
$ws.Range("E103").Value = "' 999 "
$ws.Range("E103").Style = "Normal"
$ws.Range("F103").Value = "' 9999 "
$ws.Range("F103").Style = "Normal"
$ws.Range("G103").Value = "' Clean Room Manufacturing, R&D, HR "
$ws.Range("G103").Style = "Normal"
$ws.Range("H103").Value = "' Sun_Feb_25_14:13:08_EST_2018"
$ws.Range("H103").Style = "Normal"

# Row 104
$ws.Range("A104").Value = 1519588508
$ws.Range("B104").Value = 4
$ws.Range("C104").Value = "' Indooroopilly_Shopping_Centre "
$ws.Range("C104").Style = "Normal"
$ws.Range("D104").Value = "' - "
$ws.Range("D104").Style = "Normal"
$ws.Range("E104").Value = "' 999 "
$ws.Range("E104").Style = "Normal"
$ws.Range("F104").Value = "' 9999 "
$ws.Range("F104").Style = "Normal"
$ws.Range("G104").Value = "'  "
$ws.Range("G104").Style = "Normal"
$ws.Range("H104").Value = "' Sun_Feb_25_14:55:08_EST_2018"
$ws.Range("H104").Style = "Normal"

# Row 105
$ws.Range("A105").Value = 1519588688
$ws.Range("B105").Value = 4
$ws.Range("C105").Value = "' Indooroopilly_Shopping_Centre "
$ws.Range("C105").Style = "Normal"
$ws.Range("D105").Value = "' - "
$ws.Range("D105").Style = "Normal"
$ws.Range("E105").Value = "' 999 "
$ws.Range("E105").Style = "Normal"
$ws.Range("F105").Value = "' 9999 "
$ws.Range("F105").Style = "Normal"
$ws.Range("G105").Value = "'  "
$ws.Range("G105").Style = "Normal"
$ws.Range("H105").Value = "' Sun_Feb_25_14:58:08_EST_2018"
$ws.Range("H105").Style = "Normal"

# Row 106
$ws.Range("A106").Value = 1519588701
$ws.Range("B106").Value = 4
$ws.Range("C106").Value = "' Indooroopilly_Shopping_Centre "
$ws.Range("C106").Style = "Normal"
$ws.Range("D106").Value = "' - "
$ws.Range("D106").Style = "Normal"
$ws.Range("E106").Value = "' 999 "
$ws.Range("E106").Style = "Normal"
$ws.Range("F106").Value = "' 9999 "
$ws.Range("F106").Style = "Normal"
$ws.Range("G106").Value = "'  "
$ws.Range("G106").Style = "Normal"
$ws.Range("H106").Value = "' Sun_Feb_25_14:58:21_EST_2018"
$ws.Range("H106").Style = "Normal"

# Row 107
$ws.Range("A107").Value = 1519589302
$ws.Range("B107").Value = 4
$ws.Range("C107").Value = "' Indooroopilly_Shopping_Centre "
$ws.Range("C107").Style = "Normal"
$ws.Range("D107").Value = "' - "
$ws.Range("D107").Style = "Normal"
$ws.Range("E107").Value = "' 999 "
$ws.Range("E107").Style = "Normal"
$ws.Range("F107").Value = "' 9999 "
$ws.Range("F107").Style = "Normal"
$ws.Range("G107").Value = "' Clean Room Manufacturing, Shipping, Inventory Management "
$ws.Range("G107").Style = "Normal"
$ws.Range("H107").Value = "' Sun_Feb_25_15:08:22_EST_2018"
$ws.Range("H107").Style = "Normal"

# Row 108
$ws.Range("A108").Value = 1519655158
$ws.Range("B108").Value = 4
$ws.Range("C108").Value = "' Pacific_Fair "
$ws.Range("C108").Style = "Normal"
$ws.Range("D108").Value = "' - "
$ws.Range("D108").Style = "Normal"
$ws.Range("E108").Value = "' 999 "
$ws.Range("E108").Style = "Normal"
$ws.Range("F108").Value = "' 999 "
$ws.Range("F108").Style = "Normal"
$ws.Range("G108").Value = "' Clean Room Manufacturing, Shipping, R&D "
$ws.Range("G108").Style = "Normal"
$ws.Range("H108").Value = "' Mon_Feb_26_09:25:58_EST_2018"
$ws.Range("H108").Style = "Normal"
